$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge split runs into single runs (text content unchanged, just joined)
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "Valoración de la verificación de existencia de comentarios previos en la línea de código analizada",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Valoración de la verificación de existencia de comentarios previos en la línea de código analizada",
    2) | Out-Null

$d.Content.Find.Execute(
    "La información está en orden en el comentario de la línea del programa",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "La información está en orden en el comentario de la línea del programa",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Locate the test-case table (2nd table in the document) and add a new
#    4th column ("Criterios de Aceptación")
# ---------------------------------------------------------------------------

$t = $d.Tables(2)
$t.Columns.Add() | Out-Null

# Resize columns / table to match the new layout (values are in dxa/20 = pt)
$t.Columns(2).Width = 2232 / 20.0
$t.Columns(3).Width = 3422 / 20.0
$t.Columns(4).Width = 3247 / 20.0

$t.PreferredWidthType = 3
$t.PreferredWidth = 9576 / 20.0

# ---------------------------------------------------------------------------
# 3. Populate the new column's cells with correctly-formatted content
# ---------------------------------------------------------------------------

function Set-CellXml($cell, $innerXml) {
    $r = $cell.Range
    $r.Collapse(1) | Out-Null
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml) | Out-Null
    # InsertXML leaves the pre-existing (now empty) paragraph in place before
    # the freshly inserted one - drop it so the cell ends up with exactly the
    # one paragraph we just inserted.
    $cell.Range.Paragraphs(1).Range.Delete() | Out-Null
}

$rPrCommon = '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="es-ES"/>'
$rPrHeader = '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="es-ES"/>'

# Row 1 (header row): "Criterios de Aceptación"
$cell = $t.Cell(1, 4)
$xml = '<w:p><w:pPr><w:cnfStyle w:val="100000000000"/><w:rPr>' + $rPrHeader + '</w:rPr></w:pPr>' +
       '<w:r><w:rPr>' + $rPrHeader + '</w:rPr><w:t>Criterios de Aceptación</w:t></w:r></w:p>'
Set-CellXml $cell $xml

# Row 2 (Paso 1): empty cell, formatting only
$cell = $t.Cell(2, 4)
$xml = '<w:p><w:pPr><w:cnfStyle w:val="000000100000"/><w:rPr>' + $rPrCommon + '</w:rPr></w:pPr></w:p>'
Set-CellXml $cell $xml

# Row 3 (Paso 2): empty cell, formatting only
$cell = $t.Cell(3, 4)
$xml = '<w:p><w:pPr><w:cnfStyle w:val="000000010000"/><w:rPr>' + $rPrCommon + '</w:rPr></w:pPr></w:p>'
Set-CellXml $cell $xml

# Row 4 (Paso 3): "La información debe estar en el orden apropiado."
$cell = $t.Cell(4, 4)
$xml = '<w:p><w:pPr><w:cnfStyle w:val="000000100000"/><w:rPr>' + $rPrCommon + '</w:rPr></w:pPr>' +
       '<w:r><w:rPr>' + $rPrCommon + '</w:rPr><w:t>La información debe estar en el orden apropiado.</w:t></w:r></w:p>'
Set-CellXml $cell $xml

# ---------------------------------------------------------------------------
# 4. Final trailing paragraph (after the table): language es-ES -> es-MX
# ---------------------------------------------------------------------------

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.LanguageID = "es-MX"

Write-Host "Done"
